$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the date/time number format on the existing FECHA column
# (rows 2-4) so every cell shares the same "yyyy-mm-dd h:mm:ss" style.
$ws.Range("A2:A4").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# Append the new sale record (row 5): FECHA, PRODUCTO, CANTIDAD
$ws.Range("A5").Value = 44365.9628089942
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B5").Value = "Cerveza 2"
$ws.Range("C5").Value = 10
